$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated parameters for the seasonal transition signature ("expon" row)
$ws.Range("B14").Value = 0.05
$ws.Range("C14").Value = 0.1

# Move the active selection to B15, matching the final cursor position
$ws.Range("B15").Select()
